$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "61.240.80"
$ws.Range("E2").Value = "  +1.08%  "

# Row 3
$ws.Range("D3").Value = "2.933.97"
$ws.Range("E3").Value = "  +1.04%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.99%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.12%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.29%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.71%  "

# Row 10
$ws.Range("E10").Value = "  +0.32%  "

# Row 11
$ws.Range("E11").Value = "  -1.20%  "

# Row 12
$ws.Range("E12").Value = "  +1.40%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.84"
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("E14").Value = "  -0.66%  "

# Row 15
$ws.Range("D15").Value = "3.419.04"
$ws.Range("E15").Value = "  +1.05%  "

# Row 16
$ws.Range("D16").Value = "61.228.77"
$ws.Range("E16").Value = "  +1.09%  "

# Row 17
$ws.Range("E17").Value = "  -1.46%  "

# Row 18
$ws.Range("D18").Value = "2.932.92"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "432.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.45%  "

# Row 21
$ws.Range("E21").Value = "  +1.78%  "

# Row 22
$ws.Range("E22").Value = "  -0.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "

# Row 25
$ws.Range("E25").Value = "  +1.59%  "

# Row 26
$ws.Range("E26").Value = "  +2.06%  "

# Row 27
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.43%  "

# Row 29
$ws.Range("E29").Value = "  -0.07%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.25%  "

# Row 32
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("E33").Value = "  +1.78%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0863"
$ws.Range("E34").Value = "  +3.30%  "

# Row 35
$ws.Range("E35").Value = "  +0.49%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "

# Row 37
$ws.Range("E37").Value = "  +4.45%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.125"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.19%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.26%  "

# Row 40
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.76%  "

# Row 41
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.290"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "

# Row 42
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.99%  "

# Row 43
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "377.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.93%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0347"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.26%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.718.48"
$ws.Range("E45").Value = "  +2.41%  "

# Row 46
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.97%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.97%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.106"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.15%  "

# Row 50
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.48%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.126"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.33%  "
